$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2024-12-17 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-12-18 Wednesday", 2) | Out-Null

# Update table cell values (direct cell targeting avoids any find/replace ambiguity
# since some new values equal other cells' old values)
$t = $d.Tables.Item(1)
$t.Cell(1, 1).Range.Text = "91×35=3185"  # was 96×64=6144
$t.Cell(1, 2).Range.Text = "62×74=4588"  # was 85×41=3485
$t.Cell(1, 3).Range.Text = "68×42=2856"  # was 17×54=918
$t.Cell(1, 4).Range.Text = "67×47=3149"  # was 45×18=810
$t.Cell(1, 5).Range.Text = "26×18=468"  # was 17×59=1003
$t.Cell(5, 1).Range.Text = "24×48=1152"  # was 88×33=2904
$t.Cell(5, 2).Range.Text = "74×63=4662"  # was 30×42=1260
$t.Cell(5, 3).Range.Text = "55×59=3245"  # was 77×29=2233
$t.Cell(5, 4).Range.Text = "97×90=8730"  # was 71×62=4402
$t.Cell(5, 5).Range.Text = "79×80=6320"  # was 30×65=1950
$t.Cell(10, 1).Range.Text = "81×67=5427"  # was 96×90=8640
$t.Cell(10, 2).Range.Text = "74×43=3182"  # was 23×39=897
$t.Cell(10, 3).Range.Text = "75×57=4275"  # was 58×40=2320
$t.Cell(10, 4).Range.Text = "82×30=2460"  # was 14×31=434
$t.Cell(10, 5).Range.Text = "27×55=1485"  # was 87×96=8352
$t.Cell(15, 1).Range.Text = "99×54=5346"  # was 83×19=1577
$t.Cell(15, 2).Range.Text = "98×35=3430"  # was 18×75=1350
$t.Cell(15, 3).Range.Text = "90×79=7110"  # was 12×82=984
$t.Cell(15, 4).Range.Text = "30×27=810"  # was 79×75=5925
$t.Cell(15, 5).Range.Text = "69×38=2622"  # was 44×21=924
$t.Cell(20, 1).Range.Text = "15×39=585"  # was 26×68=1768
$t.Cell(20, 2).Range.Text = "92×92=8464"  # was 88×84=7392
$t.Cell(20, 3).Range.Text = "40×13=520"  # was 66×22=1452
$t.Cell(20, 4).Range.Text = "33×32=1056"  # was 91×35=3185
$t.Cell(20, 5).Range.Text = "29×51=1479"  # was 74×45=3330
